# Make the problem statement bold-faced text.
#
# The document currently contains a single (empty) paragraph that only
# holds the "_GoBack" bookmark. The commit adds two new paragraphs in
# front of it:
#   1. a blank paragraph
#   2. a bold paragraph containing the riddle/problem statement
#
# We build the new paragraphs from raw OOXML and insert them via
# Range.InsertXML at the very start of the document so that the
# pre-existing (untouched) paragraph and its bookmark are left exactly
# as they were.

$d = $word.ActiveDocument

$problemText = "A man finds himself on a riverbank with a cat, a parrot and a bag of seed.  He needs to transport all three to the other side of the river in his boat.  However, the boat has room for only the man himself and one other item (either the cat, parrot or seed).  In his absence, the cat could eat the parrot, and the parrot would eat the bag of seed.  Show how he can get all the passengers to the other side without leaving the wrong ones alone together."

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParasXml = '<w:p ' + $wNs + '/>' + `
  '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>' + $problemText + '</w:t></w:r>' + `
  '</w:p>'

$startRange = $d.Range(0, 0)
$startRange.InsertXML($newParasXml)
